$wb = $excel.ActiveWorkbook

# --- Hoja1!A1: update the "Conversión del día" summary text ---
$hoja1 = $wb.Worksheets.Item("Hoja1")
$text = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 1.77 = 6461.31 pesos`n✅ 6461.31 pesos = 1.75 = 923.29 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"
$hoja1.Range("A1").Value = $text

# --- tasas!N10/O10/N12/O12: updated rate figures ---
$tasas = $wb.Worksheets.Item("tasas")
$tasas.Range("N10").Value = 566.5
$tasas.Range("O10").Value = 3660.33
$tasas.Range("N12").Value = 3695
$tasas.Range("O12").Value = 528
